$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number but must remain
# stored as text (matching the original inline-string cell content).
# Setting NumberFormat to "@" (Text) first stops Excel's automatic
# number coercion when the value is assigned below.
$textCells = @("D5", "D6", "D8", "D10", "D12", "D19", "D20", "D21", "D22", "D24", "D25", "D28", "D30", "D33", "D34", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, in sheet order.
$ws.Range("D2").Value = '56.519.04'
$ws.Range("E2").Value = '  -3.94%  '
$ws.Range("D3").Value = '2.355.42'
$ws.Range("E3").Value = '  -5.42%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '513.18'
$ws.Range("E5").Value = '  -4.03%  '
$ws.Range("D6").Value = '127.62'
$ws.Range("E6").Value = '  -5.71%  '
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").Value = '0.552'
$ws.Range("E8").Value = '  -2.37%  '
$ws.Range("D9").Value = '2.371.45'
$ws.Range("E9").Value = '  -5.74%  '
$ws.Range("D10").Value = '0.0958'
$ws.Range("E10").Value = '  -3.78%  '
$ws.Range("E11").Value = '  -1.34%  '
$ws.Range("D12").Value = '4.80'
$ws.Range("E12").Value = '  -8.37%  '
$ws.Range("E13").Value = '  -5.52%  '
$ws.Range("D14").Value = '2.776.37'
$ws.Range("E14").Value = '  -5.50%  '
$ws.Range("D15").Value = '56.429.69'
$ws.Range("E15").Value = '  -4.07%  '
$ws.Range("E16").Value = '  -4.55%  '
$ws.Range("E17").Value = '  -4.34%  '
$ws.Range("D18").Value = '2.359.75'
$ws.Range("E18").Value = '  -5.71%  '
$ws.Range("D19").Value = '10.24'
$ws.Range("E19").Value = '  -4.04%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '4.04'
$ws.Range("E20").Value = '  -4.80%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '309.92'
$ws.Range("E21").Value = '  -3.79%  '
$ws.Range("D22").Value = '6.08'
$ws.Range("E22").Value = '  -0.95%  '
$ws.Range("D24").Value = '64.71'
$ws.Range("E24").Value = '  -1.81%  '
$ws.Range("D25").Value = '0.998'
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("E26").Value = '  -4.90%  '
$ws.Range("D27").Value = '2.467.19'
$ws.Range("E27").Value = '  -5.90%  '
$ws.Range("D28").Value = '0.154'
$ws.Range("E28").Value = '  -4.27%  '
$ws.Range("E29").Value = '  -4.48%  '
$ws.Range("D30").Value = '174.29'
$ws.Range("E30").Value = '  +1.57%  '
$ws.Range("E31").Value = '  -5.02%  '
$ws.Range("D32").Value = '0.0₃0715'
$ws.Range("E32").Value = '  -6.81%  '
$ws.Range("D33").Value = '6.13'
$ws.Range("E33").Value = '  -3.35%  '
$ws.Range("D34").Value = '1.12'
$ws.Range("E34").Value = '  -7.31%  '
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("E36").Value = '  +0.13%  '
$ws.Range("D37").Value = '17.62'
$ws.Range("E37").Value = '  -3.20%  '
$ws.Range("D38").Value = '1.18'
$ws.Range("E38").Value = '  -6.09%  '
$ws.Range("E39").Value = '  -6.74%  '
$ws.Range("D40").Value = '0.802'
$ws.Range("E40").Value = '  +1.72%  '
$ws.Range("D41").Value = '35.46'
$ws.Range("E41").Value = '  -3.42%  '
$ws.Range("D42").Value = '1.42'
$ws.Range("E42").Value = '  -6.71%  '
$ws.Range("D43").Value = '3.34'
$ws.Range("E43").Value = '  -4.44%  '
$ws.Range("D44").Value = '4.88'
$ws.Range("E44").Value = '  -4.18%  '
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").Value = '253.40'
$ws.Range("E45").Value = '  -9.69%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '0.568'
$ws.Range("E46").Value = '  -4.30%  '
$ws.Range("D47").Value = '121.27'
$ws.Range("E47").Value = '  -8.28%  '
$ws.Range("D49").Value = '0.0487'
$ws.Range("E49").Value = '  -4.28%  '
$ws.Range("D50").Value = '0.0207'
$ws.Range("E50").Value = '  -5.36%  '
$ws.Range("D51").Value = '16.60'
$ws.Range("E51").Value = '  -6.54%  '
